$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "ECs"
$ws.Range("C2").Value = "MuSCs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02770266666666667
$ws.Range("H2").Value = 0.083108
$ws.Range("I2").Value = 0.05180130905700151
$ws.Range("J2").Value = 0.05180130905700151
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 29.52617166666667
$ws.Range("N2").Value = 88.57851500000001
$ws.Range("O2").Value = 0.3218391660320701
$ws.Range("P2").Value = 0.3218391660320701
$ws.Range("Q2").Value = 0.8179536916244445
$ws.Range("R2").Value = 7.361583224620001
$ws.Range("S2").Value = 0.01667169010627488
$ws.Range("T2").Value = 0.01667169010627488

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "ECs"
$ws.Range("C3").Value = "MuSCs"
$ws.Range("D3").Value = "Bmp15"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02770266666666667
$ws.Range("H3").Value = 0.083108
$ws.Range("I3").Value = 0.05180130905700151
$ws.Range("J3").Value = 0.05180130905700151
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 39.715023
$ws.Range("N3").Value = 119.145069
$ws.Range("O3").Value = 0.4328989896002822
$ws.Range("P3").Value = 0.4328989896002822
$ws.Range("Q3").Value = 1.100212043828
$ws.Range("R3").Value = 9.901908394452
$ws.Range("S3").Value = 0.0224247343507479
$ws.Range("T3").Value = 0.0224247343507479

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "ECs"
$ws.Range("C4").Value = "MuSCs"
$ws.Range("D4").Value = "Bmpr2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02770266666666667
$ws.Range("H4").Value = 0.083108
$ws.Range("I4").Value = 0.05180130905700151
$ws.Range("J4").Value = 0.05180130905700151
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.50081433333333
$ws.Range("N4").Value = 67.502443
$ws.Range("O4").Value = 0.2452618443676477
$ws.Range("P4").Value = 0.2452618443676476
$ws.Range("Q4").Value = 0.623332559204889
$ws.Range("R4").Value = 5.609993032844
$ws.Range("S4").Value = 0.01270488459997872
$ws.Range("T4").Value = 0.01270488459997872

$ws.Range("A5").Value = "Bmp15"
$ws.Range("B5").Value = "ECs"
$ws.Range("C5").Value = "MuSCs"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.4126986666666667
$ws.Range("H5").Value = 1.238096
$ws.Range("I5").Value = 0.7717066171516261
$ws.Range("J5").Value = 0.7717066171516261
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 29.52617166666667
$ws.Range("N5").Value = 88.57851500000001
$ws.Range("O5").Value = 0.3218391660320701
$ws.Range("P5").Value = 0.3218391660320701
$ws.Range("Q5").Value = 12.18541167860445
$ws.Range("R5").Value = 109.66870510744
$ws.Range("S5").Value = 0.2483654140855094
$ws.Range("T5").Value = 0.2483654140855094

$ws.Range("A6").Value = "Bmp15"
$ws.Range("B6").Value = "ECs"
$ws.Range("C6").Value = "MuSCs"
$ws.Range("D6").Value = "Bmp15"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.4126986666666667
$ws.Range("H6").Value = 1.238096
$ws.Range("I6").Value = 0.7717066171516261
$ws.Range("J6").Value = 0.7717066171516261
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 39.715023
$ws.Range("N6").Value = 119.145069
$ws.Range("O6").Value = 0.4328989896002822
$ws.Range("P6").Value = 0.4328989896002822
$ws.Range("Q6").Value = 16.390337038736
$ws.Range("R6").Value = 147.513033348624
$ws.Range("S6").Value = 0.3340710148327908
$ws.Range("T6").Value = 0.3340710148327908

$ws.Range("A7").Value = "Bmp15"
$ws.Range("B7").Value = "ECs"
$ws.Range("C7").Value = "MuSCs"
$ws.Range("D7").Value = "Bmpr2"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.4126986666666667
$ws.Range("H7").Value = 1.238096
$ws.Range("I7").Value = 0.7717066171516261
$ws.Range("J7").Value = 0.7717066171516261
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.50081433333333
$ws.Range("N7").Value = 67.502443
$ws.Range("O7").Value = 0.2452618443676477
$ws.Range("P7").Value = 0.2452618443676476
$ws.Range("Q7").Value = 9.286056074280891
$ws.Range("R7").Value = 83.57450466852801
$ws.Range("S7").Value = 0.189270188233326
$ws.Range("T7").Value = 0.189270188233326

$ws.Range("A8").Value = "Bmpr2"
$ws.Range("B8").Value = "ECs"
$ws.Range("C8").Value = "MuSCs"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.09438566666666666
$ws.Range("H8").Value = 0.283157
$ws.Range("I8").Value = 0.1764920737913724
$ws.Range("J8").Value = 0.1764920737913724
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 29.52617166666667
$ws.Range("N8").Value = 88.57851500000001
$ws.Range("O8").Value = 0.3218391660320701
$ws.Range("P8").Value = 0.3218391660320701
$ws.Range("Q8").Value = 2.786847396872778
$ws.Range("R8").Value = 25.081626571855
$ws.Range("S8").Value = 0.05680206184028586
$ws.Range("T8").Value = 0.05680206184028586

$ws.Range("A9").Value = "Bmpr2"
$ws.Range("B9").Value = "ECs"
$ws.Range("C9").Value = "MuSCs"
$ws.Range("D9").Value = "Bmp15"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.09438566666666666
$ws.Range("H9").Value = 0.283157
$ws.Range("I9").Value = 0.1764920737913724
$ws.Range("J9").Value = 0.1764920737913724
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 39.715023
$ws.Range("N9").Value = 119.145069
$ws.Range("O9").Value = 0.4328989896002822
$ws.Range("P9").Value = 0.4328989896002822
$ws.Range("Q9").Value = 3.748528922537
$ws.Range("R9").Value = 33.736760302833
$ws.Range("S9").Value = 0.07640324041674355
$ws.Range("T9").Value = 0.07640324041674355

$ws.Range("A10").Value = "Bmpr2"
$ws.Range("B10").Value = "ECs"
$ws.Range("C10").Value = "MuSCs"
$ws.Range("D10").Value = "Bmpr2"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.09438566666666666
$ws.Range("H10").Value = 0.283157
$ws.Range("I10").Value = 0.1764920737913724
$ws.Range("J10").Value = 0.1764920737913724
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 22.50081433333333
$ws.Range("N10").Value = 67.502443
$ws.Range("O10").Value = 0.2452618443676477
$ws.Range("P10").Value = 0.2452618443676476
$ws.Range("Q10").Value = 2.123754361394556
$ws.Range("R10").Value = 19.113789252551
$ws.Range("S10").Value = 0.04328677153434296
$ws.Range("T10").Value = 0.04328677153434295

